$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Some "Price" values are plain decimals (e.g. "299.84") that Excel would
# otherwise auto-convert to numbers, so those cells are forced to Text
# format first to keep them as literal strings like the rest of the sheet.
$ws.Range("D2").Value = '46.830.66'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '2.271.57'
$ws.Range("E3").Value = '  -3.48%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.84'
$ws.Range("E5").Value = '  -2.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.87'
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -5.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.05'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.05'
$ws.Range("E12").Value = '  -5.54%  '
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("D14").Value = '2.619.15'
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("D15").Value = '2.277.68'
$ws.Range("E15").Value = '  -2.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.66'
$ws.Range("E16").Value = '  -4.69%  '
$ws.Range("D17").Value = '46.772.75'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.799'
$ws.Range("E18").Value = '  -4.29%  '
$ws.Range("D19").Value = '0.0₃0984'
$ws.Range("E19").Value = '  +3.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.54'
$ws.Range("E20").Value = '  -7.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.83'
$ws.Range("E21").Value = '  -6.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.82'
$ws.Range("E22").Value = '  -1.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.63'
$ws.Range("E23").Value = '  +0.29%  '
$ws.Range("E24").Value = '  -6.56%  '
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -6.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.40'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.58'
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("E31").Value = '  +7.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  +7.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.54'
$ws.Range("E33").Value = '  -3.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.34'
$ws.Range("E34").Value = '  -7.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0767'
$ws.Range("E35").Value = '  -6.40%  '
$ws.Range("E36").Value = '  +2.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.64'
$ws.Range("E38").Value = '  +11.37%  '
$ws.Range("E39").Value = '  -10.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.85'
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("E41").Value = '  -7.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.08'
$ws.Range("E42").Value = '  -10.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.36'
$ws.Range("E44").Value = '  +14.78%  '
$ws.Range("D45").Value = '1.781.92'
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.89'
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '70.92'
$ws.Range("E47").Value = '  -3.79%  '
$ws.Range("E48").Value = '  -7.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.79'
$ws.Range("E49").Value = '  -2.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.58'
$ws.Range("E50").Value = '  -4.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.87'
$ws.Range("E51").Value = '  -2.02%  '
